$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "purpose" column (E2:E18) from "S.GISH" to the new value "fullRNASEQ"
$ws.Range("E2:E18").Value = "fullRNASEQ"

# Scroll/select as in the saved view
$excel.ActiveWindow.ScrollRow = 8
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("D19:F25").Select()
